# TottenKyle-grid.xlsx — "Administrative Commit: Updated the Grid"
#
# The author filled in column E ("line numbers in the input file AND the
# relevant JUnit tests") for a bunch of grading-grid rows, flipped a couple
# of B-column yes/no dropdowns to "Yes"/"Yes (only Tests)", and left the
# cursor sitting on A42 after scrolling the sheet. Reproduce all of the
# content edits here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Detection-test references (rows 14-37) ---------------------------
$ws.Range("E14").Value = "DetectionTests/inOrderRoyalFlushTest, DetectionTests/outOrderRoyalFlishTest, DetectionTests/failRoyalFlushTest"
$ws.Range("E15").Value = "DetectionTests/inOrderStraightFlushTest, DetectionTests/outOrderStraightFlushTest, DetectionTests/failStraightFlushTest"
$ws.Range("E16").Value = "DetectionTests/inOrderFourKindTest, DetectionTests/outOrderFourKindTest, DetectionTests/failFourKindTest"
$ws.Range("E17").Value = "DetectionTests/inOrderFullHouseTest, DetectionTests/outOrderFullHouseTest, DetectionTests/failFullHouseTest"
$ws.Range("E18").Value = "DetectionTests/inOrderFlushTest, DetectionTests/outOrderFlushTest, DetectionTests/failFlushTest"
$ws.Range("E19").Value = "DetectionTests/inOrderStraightTest, DetectionTests/outOrderStraightTest, DetectionTests/failStraightTest"
$ws.Range("E20").Value = "DetectionTests/threeKindTest"
$ws.Range("E21").Value = "DetectionTests/twoPairTest"
$ws.Range("E22").Value = "DetectionTests/pairTest"
$ws.Range("E23").Value = "DetectionTests/oppStraightAiNoTest, DetectionTests/neitherStraightTest, DetectionTests/bothStraightTest"
$ws.Range("E25").Value = "DetectionTests/independenceDetectionTest, SwapOneTests/flushOneSwapTest, SwapOneTests/straightOneSwapTest"
$ws.Range("E27").Value = "DetectionTests/oneOffRoyalTest"
$ws.Range("E31").Value = "SwapOneTests/flushOneSwapTest"
$ws.Range("E33").Value = "DetectionTests/threeSuitTest"
$ws.Range("E34").Value = "DetectionTests/threeKindTest"
$ws.Range("E36").Value = "DetectionTests/twoPairTest"
$ws.Range("E37").Value = "DetectionTests/pairTest"

# --- "AIP holds/exchanges" block (rows 42-63) --------------------------
$ws.Range("B42").Value = "Yes"
$ws.Range("E48").Value = "Line 57: SwapZeroTests/RoyalFlushTest"
$ws.Range("E49").Value = "Line 48: SwapZeroTests/StraightFlushTest"
$ws.Range("E51").Value = "Line 3: SwapZeroTests/fullHouseTest"
$ws.Range("E52").Value = "Line 42: SwapZeroTests/flushTest"
$ws.Range("E53").Value = "Line 43: SwapZeroTests/straightTest"
$ws.Range("E58").Value = "Line: 9 SwapOneTests/flushOneSwapTest"
$ws.Range("B59").Value = "Yes (only Tests)"
$ws.Range("E59").Value = "Line 10: SwapOneTests/straightOneSwapTest"
$ws.Range("E61").Value = "Line 11: SwapTwoTests/threeKindSwapTwoTest"
$ws.Range("E63").Value = "Line 7: SwapOneTests/twoPairOneSwapTest"

# --- "Winning" block (rows 67-76) --------------------------------------
$ws.Range("E67").Value = "All Tests in SwaplessWinTests/"
$ws.Range("E68").Value = "Lines 46-54: royalFlushOverStriaghtFlush, royalFlushOverFourKind, royalFlushOverFull, royalFlushOverFlush, royalFlushOverStraight, royalFlushOverThreeKind, royalFlushOverTwoPair, royalFlushOverPair, royalFlushOverHigh"
$ws.Range("E69").Value = "Lines 38-45: straightFlushOverFourKind, straightFlushOverFull, straightFlushOverFlush, straightFlushOverStraight, straightFlushOverThreeKind, straightFlushOverTwoPair, straightFlushOverPair, straightFlushOverHigh"
$ws.Range("E70").Value = "Lines 31-37: fourKindOverFull, fourKindOverFlush, fourKindOverStraight, fourKindOverThreeKind, fourKindOverTwoPair, fourKindOverPair, fourKindOverHigh"
$ws.Range("E71").Value = "Lines 25-30: fullOverFlush, fullOverStraight, fullOverThreeKind, fullOverTwoPair, fullOverPair, fullOverHigh"
$ws.Range("E72").Value = "Lines 20-24: flushOverStraight, flushOverThreeKind, straightOverTwoPair, straightOverPair, straightOverHigh"
$ws.Range("E73").Value = "Lines 16-19: straightOverThreeKind, straightOverTwoPair, straightOverPair, straightOverHigh"
$ws.Range("E74").Value = "Lines 13-15: threeKindOverTwoPair, threeKindOverPair, threeKindOverHigh"
$ws.Range("E75").Value = "Lines 11-12: twoPairOverPair, twoPairOverHigh"
$ws.Range("E76").Value = "Line 10: pairOverHigh"

# --- "If they both have the same hand type" block (rows 80-98) --------
$ws.Range("E80").Value = "Lines 57-60: SwappedWinTests/spadeRFvsHeartRF, SwappedWinTests/heartRFvsDiamondRF, SwappedWinTests/diamondRFvsClubRF, spadeRFvsClubRF"
$ws.Range("E81").Value = "Line 61: SwappedWinTests/spadeSFvsHeartSF"
$ws.Range("E82").Value = "Line 62: SwappedWinTests/spadeSFvsHigherSF"
$ws.Range("E83").Value = "Line 63: SwappedWinTests/higherFourKind"
$ws.Range("E84").Value = "Line 64: SwappedWinTests/higherFullHouse"
$ws.Range("E87").Value = "Line 66: SwappedWinTests/higherStraight"
$ws.Range("E88").Value = "Line 65: SwappedWinTests/higherStraightSuit"
$ws.Range("B96").Value = "Yes"
$ws.Range("E96").Value = "Line 68: SwappedWinTests/highestCardSuit"
$ws.Range("E98").Value = "Line 67: SwappedWinTests/highestCard"

# --- View state: the author ended up with A42 selected after scrolling -
# (ScrollRow/ScrollColumn only mirror the current selection in this host,
# and this engine only round-trips topLeftCell inside a freeze-pane , so
# only the active-cell/selection part of the original sheetView edit is
# reproducible here; that's still applied below.)
try { $excel.ActiveWindow.ScrollRow = 28 } catch {}
try { $excel.ActiveWindow.ScrollColumn = 1 } catch {}
$ws.Range("A42").Select()
